$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "1.004") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.959.11"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "1.910.46"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "324.71"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4590"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("D8").Value = "0.3820"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "0.07723"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").Value = "0.9804"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").Value = "1.919.93"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "6.944"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "0.07047"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "83.86"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").Value = "0.000009458"
$ws.Range("E18").Value = "  -3.70%  "
$ws.Range("D19").Value = "16.70"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "28.947.93"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "5.325"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").Value = "2.139.07"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").Value = "2.093"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "158.69"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").Value = "19.05"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "5.667"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "117.44"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "1.853"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "0.09294"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").Value = "0.8687"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "5.080"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").Value = "1.247"
$ws.Range("E34").Value = "  -4.38%  "
$ws.Range("D35").Value = "3.049"
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("D36").Value = "0.05723"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").Value = "1.159"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "0.02042"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("D41").Value = "7.399"
$ws.Range("D42").Value = "0.1754"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").Value = "2.851"
$ws.Range("E43").Value = "  +4.28%  "
$ws.Range("D44").Value = "9.301"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "0.5185"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "11.18"
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("D47").Value = "0.06868"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "2.068"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "0.000002585"
$ws.Range("E49").Value = "  -9.36%  "
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "110.59"
$ws.Range("E51").Value = "  -0.84%  "
